$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 45613
$ws.Range("C1").Value = 45620
$ws.Range("D1").Value = 45627
$ws.Range("E1").Value = 45634
$ws.Range("F1").Value = 45641
$ws.Range("G1").Value = 45648
$ws.Range("H1").Value = 45655
$ws.Range("I1").Value = 45662
$ws.Range("J1").Value = 45669
$ws.Range("K1").Value = 45676
$ws.Range("L1").Value = 45683

$ws.Range("B2").Value = 71.5
$ws.Range("C2").Value = 67.9
$ws.Range("D2").Value = 65
$ws.Range("E2").Value = 63.1
$ws.Range("F2").Value = 62.5
$ws.Range("G2").Value = 63.4
$ws.Range("H2").Value = 65.8
$ws.Range("I2").Value = 69.6
$ws.Range("J2").Value = 74.5
$ws.Range("K2").Value = 80
$ws.Range("L2").Value = 85.8

$ws.Range("B3").Value = 82.2
$ws.Range("C3").Value = 76.3
$ws.Range("D3").Value = 71.9
$ws.Range("E3").Value = 69.3
$ws.Range("F3").Value = 69.1
$ws.Range("G3").Value = 71.2
$ws.Range("H3").Value = 75.7
$ws.Range("I3").Value = 82.1
$ws.Range("J3").Value = 90
$ws.Range("K3").Value = 98.6
$ws.Range("L3").Value = 107.2

$ws.Range("B4").Value = 51.7
$ws.Range("C4").Value = 48.9
$ws.Range("D4").Value = 46.4
$ws.Range("E4").Value = 44.5
$ws.Range("F4").Value = 43.3
$ws.Range("G4").Value = 43
$ws.Range("H4").Value = 43.9
$ws.Range("I4").Value = 45.9
$ws.Range("J4").Value = 49
$ws.Range("K4").Value = 52.9
$ws.Range("L4").Value = 57.5

$ws.Range("B5").Value = 65
$ws.Range("C5").Value = 61.8
$ws.Range("D5").Value = 59.5
$ws.Range("E5").Value = 58.5
$ws.Range("F5").Value = 58.8
$ws.Range("G5").Value = 60.7
$ws.Range("H5").Value = 64.1
$ws.Range("I5").Value = 68.7
$ws.Range("J5").Value = 74.2
$ws.Range("K5").Value = 80.2
$ws.Range("L5").Value = 86.3

$ws.Range("B6").Value = 76.6
$ws.Range("C6").Value = 73.1
$ws.Range("D6").Value = 70.8
$ws.Range("E6").Value = 70.1
$ws.Range("F6").Value = 71.1
$ws.Range("G6").Value = 74
$ws.Range("H6").Value = 78.8
$ws.Range("I6").Value = 85
$ws.Range("J6").Value = 92.3
$ws.Range("K6").Value = 100
$ws.Range("L6").Value = 107.3

$ws.Range("B7").Value = 77.3
$ws.Range("C7").Value = 71.1
$ws.Range("D7").Value = 64.7
$ws.Range("E7").Value = 58.6
$ws.Range("F7").Value = 53.1
$ws.Range("G7").Value = 48.5
$ws.Range("H7").Value = 45.1
$ws.Range("I7").Value = 43
$ws.Range("J7").Value = 42.2
$ws.Range("K7").Value = 42.7
$ws.Range("L7").Value = 44.1

$ws.Range("B8").Value = 75.2
$ws.Range("C8").Value = 70.9
$ws.Range("D8").Value = 67.7
$ws.Range("E8").Value = 66
$ws.Range("F8").Value = 66
$ws.Range("G8").Value = 68
$ws.Range("H8").Value = 71.9
$ws.Range("I8").Value = 77.4
$ws.Range("J8").Value = 84.2
$ws.Range("K8").Value = 91.6
$ws.Range("L8").Value = 99.1

$ws.Range("B9").Value = 84.9
$ws.Range("C9").Value = 78.5
$ws.Range("D9").Value = 72.5
$ws.Range("E9").Value = 67.3
$ws.Range("F9").Value = 63.5
$ws.Range("G9").Value = 61.2
$ws.Range("H9").Value = 60.6
$ws.Range("I9").Value = 61.7
$ws.Range("J9").Value = 64.3
$ws.Range("K9").Value = 68
$ws.Range("L9").Value = 72.4

$ws.Range("B10").Value = 59.9
$ws.Range("C10").Value = 56.7
$ws.Range("D10").Value = 53.8
$ws.Range("E10").Value = 51.4
$ws.Range("F10").Value = 49.9
$ws.Range("G10").Value = 49.5
$ws.Range("H10").Value = 50.3
$ws.Range("I10").Value = 52.3
$ws.Range("J10").Value = 55.4
$ws.Range("K10").Value = 59.4
$ws.Range("L10").Value = 64

$ws.Range("B11").Value = 61.4
$ws.Range("C11").Value = 61.1
$ws.Range("D11").Value = 61.9
$ws.Range("E11").Value = 63.9
$ws.Range("F11").Value = 67.1
$ws.Range("G11").Value = 71.4
$ws.Range("H11").Value = 76.8
$ws.Range("I11").Value = 82.9
$ws.Range("J11").Value = 89.4
$ws.Range("K11").Value = 95.6
$ws.Range("L11").Value = 101.3

$ws.Range("B12").Value = 68.7
$ws.Range("C12").Value = 64.9
$ws.Range("D12").Value = 61.8
$ws.Range("E12").Value = 59.7
$ws.Range("F12").Value = 59
$ws.Range("G12").Value = 59.7
$ws.Range("H12").Value = 61.9
$ws.Range("I12").Value = 65.6
$ws.Range("J12").Value = 70.4
$ws.Range("K12").Value = 75.9
$ws.Range("L12").Value = 81.8

$ws.Range("B13").Value = 86.4
$ws.Range("C13").Value = 80.9
$ws.Range("D13").Value = 77.5
$ws.Range("E13").Value = 76.6
$ws.Range("F13").Value = 78.6
$ws.Range("G13").Value = 83.4
$ws.Range("H13").Value = 90.9
$ws.Range("I13").Value = 100.5
$ws.Range("J13").Value = 111.5
$ws.Range("K13").Value = 122.9
$ws.Range("L13").Value = 133.8

$ws.Range("B14").Value = 60.8
$ws.Range("C14").Value = 61.6
$ws.Range("D14").Value = 63
$ws.Range("E14").Value = 65.2
$ws.Range("F14").Value = 68.2
$ws.Range("G14").Value = 71.9
$ws.Range("H14").Value = 76.2
$ws.Range("I14").Value = 81
$ws.Range("J14").Value = 86
$ws.Range("K14").Value = 90.8
$ws.Range("L14").Value = 95.1

$ws.Range("B15").Value = 70
$ws.Range("C15").Value = 66
$ws.Range("D15").Value = 62.8
$ws.Range("E15").Value = 61
$ws.Range("F15").Value = 60.6
$ws.Range("G15").Value = 61.8
$ws.Range("H15").Value = 64.7
$ws.Range("I15").Value = 68.9
$ws.Range("J15").Value = 74.2
$ws.Range("K15").Value = 80.1
$ws.Range("L15").Value = 86

$ws.Range("B16").Value = 76.1
$ws.Range("C16").Value = 71.9
$ws.Range("D16").Value = 68.2
$ws.Range("E16").Value = 65.4
$ws.Range("F16").Value = 63.8
$ws.Range("G16").Value = 63.7
$ws.Range("H16").Value = 65.2
$ws.Range("I16").Value = 68.4
$ws.Range("J16").Value = 72.9
$ws.Range("K16").Value = 78.5
$ws.Range("L16").Value = 84.7

$ws.Range("B17").Value = 67.9
$ws.Range("C17").Value = 66.9
$ws.Range("D17").Value = 66.5
$ws.Range("E17").Value = 66.9
$ws.Range("F17").Value = 68.3
$ws.Range("G17").Value = 70.6
$ws.Range("H17").Value = 74
$ws.Range("I17").Value = 78.2
$ws.Range("J17").Value = 83.1
$ws.Range("K17").Value = 88.1
$ws.Range("L17").Value = 93

$ws.Range("B18").Value = 75.3
$ws.Range("C18").Value = 72.1
$ws.Range("D18").Value = 69
$ws.Range("E18").Value = 66.4
$ws.Range("F18").Value = 64.7
$ws.Range("G18").Value = 64.1
$ws.Range("H18").Value = 64.8
$ws.Range("I18").Value = 67.1
$ws.Range("J18").Value = 70.9
$ws.Range("K18").Value = 75.9
$ws.Range("L18").Value = 81.7

$ws.Range("B19").Value = 75.5
$ws.Range("C19").Value = 70.9
$ws.Range("D19").Value = 66.3
$ws.Range("E19").Value = 61.9
$ws.Range("F19").Value = 58.2
$ws.Range("G19").Value = 55.6
$ws.Range("H19").Value = 54.2
$ws.Range("I19").Value = 54.4
$ws.Range("J19").Value = 56.2
$ws.Range("K19").Value = 59.4
$ws.Range("L19").Value = 63.7

